$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed coin price/volume/hour figures from the latest scrape.
# Cells store these as text (matching the source sheet's inline-string layout),
# so number format is pinned to Text before writing each value to avoid
# Excel auto-coercing the numeric-looking strings into numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "307.76"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-4.57%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "9"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.09"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-6.72%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "9"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.133"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.98%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "9"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07751"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-5.45%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "9"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.249"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.90%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "9"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.624"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-11.12%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "9"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8807"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-5.57%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "9"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1029"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-7.73%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "9"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1751"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-6.43%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "9"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08978"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-5.72%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "9"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04397"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "9"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1056"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.15%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "9"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001252"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-3.56%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "9"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005837"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.25%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "9"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.354"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.29%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "9"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-4.89%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "9"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3324"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.47%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "9"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.020"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-5.27%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "9"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1339"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-3.53%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "9"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2786"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "11.80%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "9"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04162"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.34%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "9"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-3.49%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "9"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004082"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-6.03%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "9"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "8.43%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "9"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.11%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "9"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "9"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "9"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "9"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "9"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "9"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "9"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "9"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "9"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "9"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "9"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "9"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02366"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-14.37%"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "9"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05210"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-6.79%"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "9"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007937"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-4.19%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "9"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-4.93%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "9"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006343"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-2.88%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "9"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001971"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-5.72%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "9"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008351"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "10.69%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "9"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.76%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "9"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006530"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-6.62%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "9"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.07%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "9"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "98.53%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "9"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002210"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-36.45%"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "9"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.07%"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "9"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.07%"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "9"
